# "Generate Report for Handoff"
#
# A new handoff cycle was generated for the markdown source file, replacing
# the previous GUID-named source (a3da39c9-...) with a new one
# (ea06f795-...). The new handoff's xliff files use a new content hash
# (721191f8...). Because the handback step has not happened yet for this
# new cycle, the "Latest Target File" / "Latest Handback File" columns are
# cleared and "Latest Handback DateTime" resets to the zero/default
# datetime (0001-01-01 00:00:00) on both locale sheets.

$wb = $excel.ActiveWorkbook

$newGuid = "ea06f795-c856-48b4-9e4b-62cb483bb412"
$newHash = "721191f8c19ab69443f4128677d35162e543f711"

$generateDate = "2016-08-16 06:53:20"
$zhHandoffDate = "2016-08-16 06:53:15"
$deHandoffDate = "2016-08-16 06:53:20"
$resetDate = "0001-01-01 00:00:00"

# The source-file hyperlinks all still point at the same commit/URL as
# before - only the visible display text changes to the new file name.
$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfb6f57c02da00b863bd9db104973b3c50cf4fc0/e2e/$newGuid.md"

# ---------------------------------------------------------------------
# Overview sheet: file name / path columns + the handoff generation date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Reset hyperlinks on this sheet so we can rewrite B2's display text.
$overview.Range("B2").Hyperlinks.Delete()

$overview.Range("A2").Value = "$newGuid.md"
$overview.Range("B2").Value = "e2e\$newGuid.md"
$overview.Range("G2").Value = $generateDate

$overview.Hyperlinks.Add($overview.Range("B2"), $sourceUrl, "", "", "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Reset hyperlinks on this sheet (A2 source link + I2 target-file link)
# so the stale I2 link can be dropped and A2's display text rewritten.
$zhcn.Range("A2").Hyperlinks.Delete()

$zhcn.Range("A2").Value = "$newGuid.md"
$zhcn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$zhcn.Range("H2").Value = $zhHandoffDate
$zhcn.Range("I2").Value = ""
$zhcn.Range("I2").ClearFormats()
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = $resetDate

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $sourceUrl, "", "", "$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Hyperlinks.Delete()

$dede.Range("A2").Value = "$newGuid.md"
$dede.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$dede.Range("H2").Value = $deHandoffDate
$dede.Range("I2").Value = ""
$dede.Range("I2").ClearFormats()
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = $resetDate

$dede.Hyperlinks.Add($dede.Range("A2"), $sourceUrl, "", "", "$newGuid.md") | Out-Null
